$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "feature list" comment block in column A, rows 19-38.
# These are '#'-prefixed comment rows documenting the per-dataset feature
# names used elsewhere in the sheet (review-response follow-up).
$ws.Range("A19").Value = '#   WHONDRS features'
$ws.Range("A20").Value = '#       ''DS_Water.Column.Height_cm'''
$ws.Range("A21").Value = '#       ''SW_Temp_degC''            '
$ws.Range("A22").Value = '#'
$ws.Range("A23").Value = '#   StreamStats features'
$ws.Range("A24").Value = '#       ''lat''     '
$ws.Range("A25").Value = '#       ''DRNAREA'''
$ws.Range("A26").Value = '#'
$ws.Range("A27").Value = '#   HydroSheds features'
$ws.Range("A28").Value = '#       ''Natural discharge - dis_m3_pmn'''
$ws.Range("A29").Value = '#       ''Air temp - tmp_dc_s05''         '
$ws.Range("A30").Value = '#       ''Potential ET - pet_mm_s04''      '
$ws.Range("A31").Value = '#       ''Actual ET - aet_mm_s04''        '
$ws.Range("A32").Value = '#'
$ws.Range("A33").Value = '#   EPAW-Catchment features'
$ws.Range("A34").Value = '#       ''Mean summer stream temperature - MSST'''
$ws.Range("A35").Value = '#'
$ws.Range("A36").Value = '#   EPAW-Watershed features'
$ws.Range("A37").Value = '#       ''PctImpWs'''
$ws.Range("A38").Value = '#       ''TmeanWs'''

# Move the active selection, as captured when the workbook was last saved.
[void]$ws.Range("C14").Select()
